$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2579443333333333
$ws.Range("H2").Value = 0.773833
$ws.Range("I2").Value = 0.05467096027587241
$ws.Range("J2").Value = 0.0546709602758724
$ws.Range("M2").Value = 35.32109533333334
$ws.Range("N2").Value = 105.963286
$ws.Range("O2").Value = 0.4123245624288747
$ws.Range("P2").Value = 0.4123245624288747
$ws.Range("Q2").Value = 9.110876388359779
$ws.Range("R2").Value = 81.997887495238
$ws.Range("S2").Value = 0.02254217977331548
$ws.Range("T2").Value = 0.02254217977331548
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2579443333333333
$ws.Range("H3").Value = 0.773833
$ws.Range("I3").Value = 0.05467096027587241
$ws.Range("J3").Value = 0.0546709602758724
$ws.Range("O3").Value = 0.01433703690686912
$ws.Range("P3").Value = 0.01433703690686912
$ws.Range("Q3").Value = 0.3167964825194444
$ws.Range("R3").Value = 2.851168342675
$ws.Range("S3").Value = 0.0007838195752091582
$ws.Range("T3").Value = 0.0007838195752091581
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2579443333333333
$ws.Range("H4").Value = 0.773833
$ws.Range("I4").Value = 0.05467096027587241
$ws.Range("J4").Value = 0.0546709602758724
$ws.Range("M4").Value = 2.583168333333334
$ws.Range("N4").Value = 7.749505
$ws.Range("O4").Value = 0.03015489023401347
$ws.Range("P4").Value = 0.03015489023401347
$ws.Range("Q4").Value = 0.6663136336294445
$ws.Range("R4").Value = 5.996822702665
$ws.Range("S4").Value = 0.001648596806107043
$ws.Range("T4").Value = 0.001648596806107043
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2579443333333333
$ws.Range("H5").Value = 0.773833
$ws.Range("I5").Value = 0.05467096027587241
$ws.Range("J5").Value = 0.0546709602758724
$ws.Range("M5").Value = 46.53090866666667
$ws.Range("N5").Value = 139.592726
$ws.Range("O5").Value = 0.5431835104302428
$ws.Range("P5").Value = 0.5431835104302427
$ws.Range("Q5").Value = 12.00238421541756
$ws.Range("R5").Value = 108.021457938758
$ws.Range("S5").Value = 0.02969636412124073
$ws.Range("T5").Value = 0.02969636412124072
$ws.Range("I6").Value = 0.4403936734732808
$ws.Range("J6").Value = 0.4403936734732807
$ws.Range("M6").Value = 35.32109533333334
$ws.Range("N6").Value = 105.963286
$ws.Range("O6").Value = 0.4123245624288747
$ws.Range("P6").Value = 0.4123245624288747
$ws.Range("Q6").Value = 73.39129038495224
$ws.Range("R6").Value = 660.5216134645701
$ws.Range("S6").Value = 0.1815851287113152
$ws.Range("T6").Value = 0.1815851287113152
$ws.Range("I7").Value = 0.4403936734732808
$ws.Range("J7").Value = 0.4403936734732807
$ws.Range("O7").Value = 0.01433703690686912
$ws.Range("P7").Value = 0.01433703690686912
$ws.Range("S7").Value = 0.006313940350138094
$ws.Range("T7").Value = 0.006313940350138094
$ws.Range("I8").Value = 0.4403936734732808
$ws.Range("J8").Value = 0.4403936734732807
$ws.Range("M8").Value = 2.583168333333334
$ws.Range("N8").Value = 7.749505
$ws.Range("O8").Value = 0.03015489023401347
$ws.Range("P8").Value = 0.03015489023401347
$ws.Range("Q8").Value = 5.367388963330557
$ws.Range("R8").Value = 48.306500669975
$ws.Range("S8").Value = 0.01328002288334075
$ws.Range("T8").Value = 0.01328002288334075
$ws.Range("I9").Value = 0.4403936734732808
$ws.Range("J9").Value = 0.4403936734732807
$ws.Range("M9").Value = 46.53090866666667
$ws.Range("N9").Value = 139.592726
$ws.Range("O9").Value = 0.5431835104302428
$ws.Range("P9").Value = 0.5431835104302427
$ws.Range("Q9").Value = 96.68339550637447
$ws.Range("R9").Value = 870.15055955737
$ws.Range("S9").Value = 0.2392145815284868
$ws.Range("T9").Value = 0.2392145815284867
$ws.Range("G10").Value = 2.217259
$ws.Range("H10").Value = 6.651777
$ws.Range("I10").Value = 0.4699451123575263
$ws.Range("J10").Value = 0.4699451123575263
$ws.Range("M10").Value = 35.32109533333334
$ws.Range("N10").Value = 105.963286
$ws.Range("O10").Value = 0.4123245624288747
$ws.Range("P10").Value = 0.4123245624288747
$ws.Range("Q10").Value = 78.31601651769134
$ws.Range("R10").Value = 704.8441486592221
$ws.Range("S10").Value = 0.1937699128184054
$ws.Range("T10").Value = 0.1937699128184054
$ws.Range("G11").Value = 2.217259
$ws.Range("H11").Value = 6.651777
$ws.Range("I11").Value = 0.4699451123575263
$ws.Range("J11").Value = 0.4699451123575263
$ws.Range("O11").Value = 0.01433703690686912
$ws.Range("P11").Value = 0.01433703690686912
$ws.Range("Q11").Value = 2.723145118008333
$ws.Range("R11").Value = 24.508306062075
$ws.Range("S11").Value = 0.00673762042007261
$ws.Range("T11").Value = 0.00673762042007261
$ws.Range("G12").Value = 2.217259
$ws.Range("H12").Value = 6.651777
$ws.Range("I12").Value = 0.4699451123575263
$ws.Range("J12").Value = 0.4699451123575263
$ws.Range("M12").Value = 2.583168333333334
$ws.Range("N12").Value = 7.749505
$ws.Range("O12").Value = 0.03015489023401347
$ws.Range("P12").Value = 0.03015489023401347
$ws.Range("Q12").Value = 5.727553235598333
$ws.Range("R12").Value = 51.547979120385
$ws.Range("S12").Value = 0.01417114327915234
$ws.Range("T12").Value = 0.01417114327915233
$ws.Range("G13").Value = 2.217259
$ws.Range("H13").Value = 6.651777
$ws.Range("I13").Value = 0.4699451123575263
$ws.Range("J13").Value = 0.4699451123575263
$ws.Range("M13").Value = 46.53090866666667
$ws.Range("N13").Value = 139.592726
$ws.Range("O13").Value = 0.5431835104302428
$ws.Range("P13").Value = 0.5431835104302427
$ws.Range("Q13").Value = 103.1710760193447
$ws.Range("R13").Value = 928.539684174102
$ws.Range("S13").Value = 0.2552664358398961
$ws.Range("T13").Value = 0.255266435839896
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1650883333333333
$ws.Range("H14").Value = 0.495265
$ws.Range("I14").Value = 0.03499025389332058
$ws.Range("J14").Value = 0.03499025389332058
$ws.Range("M14").Value = 35.32109533333334
$ws.Range("N14").Value = 105.963286
$ws.Range("O14").Value = 0.4123245624288747
$ws.Range("P14").Value = 0.4123245624288747
$ws.Range("Q14").Value = 5.831100760087779
$ws.Range("R14").Value = 52.47990684079001
$ws.Range("S14").Value = 0.01442734112583864
$ws.Range("T14").Value = 0.01442734112583863
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1650883333333333
$ws.Range("H15").Value = 0.495265
$ws.Range("I15").Value = 0.03499025389332058
$ws.Range("J15").Value = 0.03499025389332058
$ws.Range("O15").Value = 0.01433703690686912
$ws.Range("P15").Value = 0.01433703690686912
$ws.Range("Q15").Value = 0.2027546123194444
$ws.Range("R15").Value = 1.824791510875
$ws.Range("S15").Value = 0.000501656561449258
$ws.Range("T15").Value = 0.000501656561449258
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1650883333333333
$ws.Range("H16").Value = 0.495265
$ws.Range("I16").Value = 0.03499025389332058
$ws.Range("J16").Value = 0.03499025389332058
$ws.Range("M16").Value = 2.583168333333334
$ws.Range("N16").Value = 7.749505
$ws.Range("O16").Value = 0.03015489023401347
$ws.Range("P16").Value = 0.03015489023401347
$ws.Range("Q16").Value = 0.4264509548694445
$ws.Range("R16").Value = 3.838058593825
$ws.Range("S16").Value = 0.001055127265413345
$ws.Range("T16").Value = 0.001055127265413344
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1650883333333333
$ws.Range("H17").Value = 0.495265
$ws.Range("I17").Value = 0.03499025389332058
$ws.Range("J17").Value = 0.03499025389332058
$ws.Range("M17").Value = 46.53090866666667
$ws.Range("N17").Value = 139.592726
$ws.Range("O17").Value = 0.5431835104302428
$ws.Range("P17").Value = 0.5431835104302427
$ws.Range("Q17").Value = 7.681710160265556
$ws.Range("R17").Value = 69.13539144239
$ws.Range("S17").Value = 0.01900612894061935
$ws.Range("T17").Value = 0.01900612894061934

Write-Host "Updated 190 cell values with new TPM-derived NATMI statistics."
